$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '72.939.84'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '4.043.40'
$ws.Range('E3').Value = '  +0.54%  '
$ws.Range('E4').Value = '  -0.03%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.27'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  +4.12%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '150.75'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  -1.10%  '
$ws.Range('D7').Value = '4.035.65'
$ws.Range('E7').Value = '  +0.58%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.696'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -1.64%  '
$ws.Range('E9').Value = '  -0.01%  '
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.764'
$ws.Range('D10').Style = $origStyle
$ws.Range('E10').Value = '  +1.44%  '
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.174'
$ws.Range('D11').Style = $origStyle
$ws.Range('E11').Value = '  +1.08%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.99'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +13.42%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000331'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +1.97%  '
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '10.99'
$ws.Range('D14').Style = $origStyle
$ws.Range('E14').Value = '  +2.58%  '
$ws.Range('D15').Value = '4.695.90'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '4.041.17'
$ws.Range('E16').Value = '  +0.73%  '
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.54'
$ws.Range('D17').Style = $origStyle
$ws.Range('E17').Value = '  +3.25%  '
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.79'
$ws.Range('D18').Style = $origStyle
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  +1.57%  '
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('D21').Value = '72.836.16'
$ws.Range('E21').Value = '  +1.20%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '445.67'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  +3.92%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '98.16'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -0.95%  '
$ws.Range('E24').Value = '  +0.66%  '
$ws.Range('E25').Value = '  +3.98%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.77'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  +2.24%  '
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.40'
$ws.Range('D27').Style = $origStyle
$ws.Range('E27').Value = '  +22.01%  '
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.34'
$ws.Range('D28').Style = $origStyle
$ws.Range('E28').Value = '  +2.37%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '10.98'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +1.33%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.94'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  +1.61%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '37.31'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +1.01%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.97'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  +14.52%  '
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.134'
$ws.Range('D33').Style = $origStyle
$ws.Range('E33').Value = '  +3.44%  '
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '13.64'
$ws.Range('D34').Style = $origStyle
$ws.Range('E34').Value = '  +1.38%  '
$origStyle = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '687.77'
$ws.Range('D35').Style = $origStyle
$ws.Range('E35').Value = '  +1.23%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '48.98'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +14.60%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '67.55'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('D38').Value = '0.0₃0923'
$ws.Range('E38').Value = '  +12.43%  '
$origStyle = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.451'
$ws.Range('D39').Style = $origStyle
$ws.Range('E39').Value = '  +6.50%  '
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('E41').Value = '  +1.07%  '
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.39'
$ws.Range('D42').Style = $origStyle
$ws.Range('E42').Value = '  -0.49%  '
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = $origStyle
$ws.Range('E43').Value = '  +0.22%  '
$ws.Range('E44').Value = '  +2.34%  '
$ws.Range('E45').Value = '  +14.72%  '
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('E47').Value = '  +0.15%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.71'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.14'
$ws.Range('D49').Style = $origStyle
$ws.Range('E49').Value = '  +4.74%  '
$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.55'
$ws.Range('D50').Style = $origStyle
$ws.Range('E50').Value = '  +6.49%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '3.35'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -2.03%  '
